$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update account stats after closing a trade (a loss), per "Fix trade closing functionality"
$ws.Range("B2").Value = 1201
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 66.66666666666666
